$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 4 (B4:E4)
$ws.Range("B4").Value = 1.1499999999999999
$ws.Range("C4").Value = 2.4
$ws.Range("D4").Value = 0.15
$ws.Range("E4").Value = 0.55000000000000004

# Apply a solid red fill to the H2:K5 range
$ws.Range("H2:K5").Interior.Color = 255

# Update the active selection to E9
$ws.Range("E9").Select()
